$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use Range.Copy (cell-to-cell) throughout instead of plain .Value assignment
# so that text that looks like a date/number (e.g. "01/01/2016") keeps its
# original shared-string type/style instead of being auto-coerced into a
# date serial the way a literal .Value write would be.
#
# Because several source cells are read into more than one destination (and
# some cells are simultaneously a destination for one move and a source for
# another), every copy below is ordered so the original content of a cell is
# always read out *before* that cell is overwritten or cleared.

# Row 10 (B/C): replace the long objectives paragraph with the professor
# line that currently lives at B13/C13.
$ws.Range("B13").Copy($ws.Range("B10"))
$ws.Range("C13").Copy($ws.Range("C10"))

# Row 15 (B/C) will need that same professor line later - copy it out now,
# before B13/C13 gets overwritten.
$ws.Range("B13").Copy($ws.Range("B15"))
$ws.Range("C13").Copy($ws.Range("C15"))

# Row 13 (A): old row 15's "Programa resumido:" label - copy before A15 is
# overwritten.
$ws.Range("A15").Copy($ws.Range("A13"))

# Row 13 (B/C): reuse the existing "01/01/2016" text cell (B8/C8).
$ws.Range("B8").Copy($ws.Range("B13"))
$ws.Range("C8").Copy($ws.Range("C13"))

# Row 14 (A): old row 16's "Short syllabus:" label.
$ws.Range("A16").Copy($ws.Range("A14"))

# Row 15 (A): old row 17's "Programa:" label.
$ws.Range("A17").Copy($ws.Range("A15"))

# Row 16 (A): old row 18's "Syllabus:" label.
$ws.Range("A18").Copy($ws.Range("A16"))

# Row 17 (A): old row 19's "Avaliação:" label.
$ws.Range("A19").Copy($ws.Range("A17"))

# Row 18 (A): old row 20's "Método:" label.
$ws.Range("A20").Copy($ws.Range("A18"))

# Row 18 (B/C): the Katia professor line, currently at B14/C14 - copy before
# B14/C14 get cleared.
$ws.Range("B14").Copy($ws.Range("B18"))
$ws.Range("C14").Copy($ws.Range("C18"))
$ws.Range("B14").Clear()
$ws.Range("C14").Clear()

# Row 19 (A): old row 21's "Critério:" label.
$ws.Range("A21").Copy($ws.Range("A19"))

# Row 19 (B/C): old row 20's teaching-method text.
$ws.Range("B20").Copy($ws.Range("B19"))
$ws.Range("C20").Copy($ws.Range("C19"))

# Row 20 (A): old row 22's "Norma de recuperação:" label.
$ws.Range("A22").Copy($ws.Range("A20"))

# Row 20 (B/C): old row 21's criterion text.
$ws.Range("B21").Copy($ws.Range("B20"))
$ws.Range("C21").Copy($ws.Range("C20"))

# Row 21 (A): old row 23's "Bibliografia:" label.
$ws.Range("A23").Copy($ws.Range("A21"))

# Row 21 (B/C): old row 22's recovery-rule text.
$ws.Range("B22").Copy($ws.Range("B21"))
$ws.Range("C22").Copy($ws.Range("C21"))

# Row 22 (A): old row 24's "Requisitos:" label.
$ws.Range("A24").Copy($ws.Range("A22"))

# Rows 17/22 have no B/C content in the new layout.
$ws.Range("B17").Clear()
$ws.Range("C17").Clear()
$ws.Range("B22").Clear()
$ws.Range("C22").Clear()

# Row 23 (B/C): old row 25's requirement text; row 23 has no A label.
$ws.Range("B25").Copy($ws.Range("B23"))
$ws.Range("C25").Copy($ws.Range("C23"))
$ws.Range("A23").Clear()

# Remove the now-unused trailing rows 24/25 (delete higher row first so the
# row numbers of rows above aren't disturbed mid-way).
$ws.Rows.Item(25).Delete()
$ws.Rows.Item(24).Delete()

# Row heights per the target layout. AutoFit() on a row clears any explicit
# height back to the sheet default (rows 17/22 have no customHeight in the
# target).
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(14).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(16).RowHeight = 120
$ws.Rows.Item(17).AutoFit()
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(19).RowHeight = 60
$ws.Rows.Item(20).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
$ws.Rows.Item(22).AutoFit()
$ws.Rows.Item(23).RowHeight = 30
